# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" / "Error Detail" columns (I/J/K/P) for the 081c8b3d-90af-... row
# (row 6) on both the zh-cn and de-de sheets, now that a handback xliff has
# been generated for them, and widens the Error Detail column to fit the
# message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c5065c227dc81c154615fdef90e0287d8c9ff37/e2e/081c8b3d-90af-40d8-9de9-529aa9ff6b65.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31281eceef8c628549c3bee9fba324fa6d06c19e/e2e/081c8b3d-90af-40d8-9de9-529aa9ff6b65.md."

# --- zh-cn sheet ---------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P) to fit the long message.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# I6 ("Latest Target File"): hyperlink to the handback markdown file.
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c5065c227dc81c154615fdef90e0287d8c9ff37/e2e/081c8b3d-90af-40d8-9de9-529aa9ff6b65.md", "", "", "081c8b3d-90af-40d8-9de9-529aa9ff6b65.md")

# J6 ("Latest Handback File"): generated handback xliff file name.
$wsZhCn.Range("J6").Value = "081c8b3d-90af-40d8-9de9-529aa9ff6b65.fb4af54f8e4f0eae660e18dd04e1f782a25ee71e.zh-cn.xlf"

# K6 ("Latest Handback DateTime"): when the handback xliff was generated.
$wsZhCn.Range("K6").Value = "2016-09-04 20:48:56"
$wsZhCn.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# P6 ("Error Detail"): handback file is behind the latest source.
$wsZhCn.Range("P6").Value = $errorDetail

# --- de-de sheet ----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Widen the "Error Detail" column (P) to fit the long message.
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664

# I6 ("Latest Target File"): hyperlink to the handback markdown file.
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c5065c227dc81c154615fdef90e0287d8c9ff37/e2e/081c8b3d-90af-40d8-9de9-529aa9ff6b65.md", "", "", "081c8b3d-90af-40d8-9de9-529aa9ff6b65.md")

# J6 ("Latest Handback File"): generated handback xliff file name.
$wsDeDe.Range("J6").Value = "081c8b3d-90af-40d8-9de9-529aa9ff6b65.fb4af54f8e4f0eae660e18dd04e1f782a25ee71e.de-de.xlf"

# K6 ("Latest Handback DateTime"): when the handback xliff was generated.
$wsDeDe.Range("K6").Value = "2016-09-04 20:49:07"
$wsDeDe.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# P6 ("Error Detail"): handback file is behind the latest source.
$wsDeDe.Range("P6").Value = $errorDetail

Write-Host "Localization status report updated for 081c8b3d-90af-40d8-9de9-529aa9ff6b65"
